$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Gehan Adel, Administrator, Dr. Amira Sobhy, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat'
$ws.Range("G3").Value = 'Dr. Majorelle Magdy, Dr. Eman Tantawi, Administrator, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Asmaa Reda'
$ws.Range("G4").Value = 'Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda'
$ws.Range("G5").Value = 'Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat'
$ws.Range("G6").Value = 'Dr. Majorelle Magdy, Dr. Manar Montaser, Dr. Menna tuâ€™Allah Medhat, Dr. Alshimaa Atef, Dr. Mohammad El-Tanany'
$ws.Range("G7").Value = 'Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Nada Mohammad, Dr. Fatma Elhady, Dr. Lamiaa Ossama, Dr. Menna tu''Alllah Mohammad, Dr. Kerelos Zareef'
$ws.Range("G8").Value = 'Dr. Abeer Ragab, Dr. Nada Mohammad'
$ws.Range("G9").Value = 'Dr. Shimaa Ashraf, Dr. Safa Hany'
$ws.Range("G11").Value = 'Dr. Aya Saeed, Dr. Amal Awwad, Dr. Safa Hany'
$ws.Range("G12").Value = 'Dr. Madeha Saeed, Dr. Amira Ibrahim, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel'
$ws.Range("G13").Value = 'Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim'
$ws.Range("G15").Value = 'Dr. Rania Ahmad Youssef, Dr. Mohammad Safwat'
$ws.Range("G17").Value = 'Dr. Esraa Samy, Dr. Mohammad Safwat'
$ws.Range("G19").Value = 'Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef'
$ws.Range("G20").Value = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Range("G25").Value = 'Dr. Nouran Mahmoud, Menna tuâ€™Allah Gamil'
$ws.Range("G28").Value = 'Dr. Maryam Ashraf, Dr. Aya Emad'
$ws.Range("G30").Value = 'Dr. Shorok Mohammad, Dr. Wafaa Ebida, Dr. Yassmen Ahmad, Dr. Aya Hanafy'
